$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 5: new time-registration entry "Review UC12"
$ws.Range("A5").Value = "Review UC12"
$ws.Range("B5").Value = "Reviewer"
$ws.Range("C5").Value = 43896
$ws.Range("D5").Value = 0.375
$ws.Range("E5").Value = 0.41666666666666669
$ws.Range("F5").Value = 30

# Row 6: new time-registration entry "Review + Edit DOM12"
$ws.Range("A6").Value = "Review + Edit DOM12"
$ws.Range("B6").Value = "Reviewer"
$ws.Range("C6").Value = 43896
$ws.Range("D6").Value = 0.41666666666666669
$ws.Range("E6").Value = 0.5625
$ws.Range("F6").Value = 60

# The two newly-filled rows no longer carry the column date/time
# data-validation rules (they now hold real, already-validated data),
# splitting the C/D validation ranges around rows 5:6.
$ws.Range("C5:C6").Validation.Delete()
$ws.Range("D5:D6").Validation.Delete()

# Drop the unused trailing placeholder rows (33:47) so the sheet's
# used range shrinks back down to A1:H32.
$ws.Range("C33:C47").Clear()

# Leave the selection where the author's cursor ended up.
[void]$ws.Range("D3").Select()
